$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "postStimBlankT" column (J),
# shifting it and everything after it (maskRR, maskOnOff, nRevs,
# priorMean, priorSD) two columns to the right.
$ws.Range("J1").EntireColumn.Insert()
$ws.Range("J1").EntireColumn.Insert()

# New header cells for the inserted columns.
$ws.Range("J1").Value = "stimRamp"
$ws.Range("K1").Value = "rampLin"

# New data columns: stimRamp / rampLin, both constant across all rows.
$ws.Range("J2:J5").Value = 1
$ws.Range("K2:K5").Value = 1

# jitTmax (column E) changed from 500 to 250 for every condition row.
$ws.Range("E2:E5").Value = 250

# nRevs (now column O, previously M) changed from 16 to 20.
$ws.Range("O2:O5").Value = 20

# priorMean (now column P, previously N) changed for the first two
# conditions only (rows 2 and 3), from 0.2 to 0.1.
$ws.Range("P2").Value = 0.1
$ws.Range("P3").Value = 0.1

# priorSD (now column Q, previously O) changed for the first two
# conditions only (rows 2 and 3), from 1 to 0.5.
$ws.Range("Q2").Value = 0.5
$ws.Range("Q3").Value = 0.5

# Update the saved selection to match the target workbook.
[void]$ws.Range("E2").Select()
